# Fix the "harvester" column (column B) in the rnaSamples sheet: Holly's
# samples were harvested by "S.GISH", but the column still showed the old
# "Retrofitted_1337" placeholder value (shared with the unrelated
# rnaPreparer/rnaPrepMethod columns). Update B2:B24 to the correct harvester.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$harvesterRange = $ws.Range("B2:B24")
$harvesterRange.Value = "S.GISH"

# Mirror the author's workflow of selecting the harvester column after
# making the correction.
$ws.Columns.Item(2).Select() | Out-Null
